# Add 2022-Q3 sheet/data for 600663-陆家嘴 holdings workbook.

$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item(1)   # "总计" summary sheet, stays first
$q1Sheet    = $wb.Worksheets.Item(2)   # current "2022-Q1" sheet (will end up 3rd)

# 1. Build the new "2022-Q3" worksheet while it is appended at the END of the
#    workbook (inserting it directly into slot 2 triggers a COM
#    paste/format quirk in this host), then relocate it with Move() once its
#    content is ready. Final tab order:
#    总计, 2022-Q3, 2022-Q1, 2021-Q3, 2021-Q2, 2020-Q4
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "2022-Q3"

# Copy header-row / first-column formatting from an existing detail sheet so
# the new sheet matches the look of its siblings.
$q1Sheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$q1Sheet.Range("A2:A3").Copy()
$newSheet.Range("A2:A7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 2. Populate the header row.
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# 3. Populate the data rows. Columns B (fund code, has leading zeros) and
#    D, E, F, G (percent/scale figures) are stored as text in the source
#    data, so force a text number-format before assigning them.
$newSheet.Range("B2:B7").NumberFormat = "@"
$newSheet.Range("D2:G7").NumberFormat = "@"

function Set-FundRow($row, $idx, $code, $name, $scale, $pos, $pct, $mv, $rank) {
    $newSheet.Range("A$row").Value = $idx
    $newSheet.Range("B$row").Value = $code
    $newSheet.Range("C$row").Value = $name
    $newSheet.Range("D$row").Value = $scale
    $newSheet.Range("E$row").Value = $pos
    $newSheet.Range("F$row").Value = $pct
    $newSheet.Range("G$row").Value = $mv
    $newSheet.Range("H$row").Value = $rank
}

Set-FundRow 2 0 "510810" "汇添富中证上海国企ETF"               "63.53" "98.46" "2.77" "1.7598" 10
Set-FundRow 3 1 "012250" "安信平衡增利混合A"                    "2.63"  "60.12" "2.32" "0.0610" 8
Set-FundRow 4 2 "515450" "南方标普中国A股大盘红利低波50ETF"     "2.17"  "99.66" "2.79" "0.0605" 6
Set-FundRow 5 3 "012251" "安信平衡增利混合C"                    "2.10"  "60.12" "2.32" "0.0487" 8
Set-FundRow 6 4 "009658" "汇丰晋信中小盘低波动策略股票A"        "0.85"  "90.14" "1.99" "0.0169" 2
Set-FundRow 7 5 "009775" "汇丰晋信中小盘低波动策略股票C"        "0.04"  "90.14" "1.99" "0.0008" 2

# 4. Move the finished sheet into its proper tab position: right after
#    "总计" and right before "2022-Q1".
$newSheet.Move($q1Sheet)

# 5. Update the "总计" summary sheet: insert the 2022-Q3 row after the header
#    and shift the remaining quarters down by one, adding 2020-Q4 back at the
#    new row 6.
$totalSheet.Range("A6").Value = 4
$totalSheet.Range("B6").Value = "2020-Q4"
$totalSheet.Range("C6").Value = 1
$totalSheet.Range("D6").Value = 3.04

$totalSheet.Range("A5").Value = 3
$totalSheet.Range("B5").Value = "2021-Q2"
$totalSheet.Range("C5").Value = 2
$totalSheet.Range("D5").Value = 2.51

$totalSheet.Range("A4").Value = 2
$totalSheet.Range("B4").Value = "2021-Q3"
$totalSheet.Range("C4").Value = 1
$totalSheet.Range("D4").Value = 2.4

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2022-Q1"
$totalSheet.Range("C3").Value = 2
$totalSheet.Range("D3").Value = 2.2

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q3"
$totalSheet.Range("C2").Value = 6
$totalSheet.Range("D2").Value = 1.95

# A6 needs the same style as the rest of column A in this sheet.
$totalSheet.Range("A5").Copy()
$totalSheet.Range("A6").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$totalSheet.Range("A6").Value = 4

# 6. Restore the originally-active tab. Inserting/moving a sheet above can
#    shift which physical sheet the engine considers "active" (it tracked a
#    raw tab index rather than sheet identity), so re-fetch "2020-Q4" by
#    name and re-activate it explicitly.
$wb.Worksheets.Item("2020-Q4").Activate()
